{"js": "// Load all paragraphs in the document body so we can find the\n// \"Shooter has lives\" bullet and insert a new sibling bullet after it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(p => p.text.trim() === \"Shooter has lives\");\nif (!target) {\n  throw new Error('Could not find paragraph with text \"Shooter has lives\"');\n}\n\n// Insert a new paragraph right after it; Word carries over the source\n// paragraph's list formatting (ListParagraph style, ilvl 1, numId 1)\n// to the newly inserted paragraph automatically.\ntarget.insertParagraph(\"Added Power ups\", \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Shooter has lives\" bullet using Find (collapses/updates the\n# range in place to the matched text span, just like real Word COM).\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Shooter has lives\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find paragraph 'Shooter has lives'\"\n}\n\n# Map the matched range back to its 1-based Paragraphs index.\n$targetIndex = 0\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n        $targetIndex = $i\n        break\n    }\n    $i++\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not resolve matched range to a paragraph\"\n}\n\n# Insert a new paragraph right after it; it inherits the source\n# paragraph's list formatting (ListParagraph style, ilvl 1, numId 1).\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.InsertParagraphAfter() | Out-Null\n\n# Re-fetch the freshly inserted paragraph from the collection and set its\n# text (writing through the range returned by InsertParagraphAfter\n# directly does not stick, so address it via Paragraphs.Item instead).\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Added Power ups\"\n"}
